$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink on the "SBS..." source URL cell (A52) entirely.
$ws.Hyperlinks.Delete()

# Insert a new blank row above row 51, shifting the "Source:" block down by one row.
$ws.Rows.Item(51).Insert()

# After the insert:
#   A50 = "Source:"                                            (unchanged)
#   A51 = <new blank row>
#   A52 = "SBS Main Indicators, ..."                            (was A51)
#   A53 = "http://epp.eurostat.ec.europa.eu/..."                (was A52, hyperlink style)
#   A54 = ""                                                    (was A53)
#   A57 = "SME Performance Review EU"                           (was A56)
#   A58 = "SME Performance Review EU, \"SBA Fact sheet\", ..."  (was A57)

# Move the URL text down into A54 (which already carries the plain "source" style).
$ws.Range("A54").Value = $ws.Range("A53").Value()

# Clear the old URL text out of A53 and strip the leftover hyperlink formatting
# (blue underline) so it matches the plain "source" style used elsewhere.
$ws.Range("A53").Value = ""
$ws.Range("A53").Font.Underline = $false
$ws.Range("A53").Font.Color = 0
$ws.Range("A53").Font.Italic = $true

# The final row duplicates the "SME Performance Review EU" label (instead of the
# old, longer citation sentence) directly below the existing bold header row.
$ws.Range("A58").Value = "SME Performance Review EU"
